$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1438.8552
$ws.Range("J17").Value = 1438.8552
$ws.Range("L17").Value = 4316.5656
$ws.Range("N17").Value = -4652.5656

# Row 39
$ws.Range("H39").Value = 532.94446
$ws.Range("I39").Value = 137.2
$ws.Range("J39").Value = 1027.625
$ws.Range("K39").Value = 411.6
$ws.Range("L39").Value = 3082.875
$ws.Range("M39").Value = -115.6
$ws.Range("N39").Value = -3674.875

# Row 44
$ws.Range("H44").Value = 12345
$ws.Range("J44").Value = 12345
$ws.Range("L44").Value = 12345
$ws.Range("N44").Value = -13269

# Row 47
$ws.Range("H47").Value = 30030
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 30030
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 30030
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -31974

# Row 52
$ws.Range("H52").Value = 433.33334
$ws.Range("I52").Value = 433.33334
$ws.Range("K52").Value = 1300.00002
$ws.Range("M52").Value = -1140.00002

# Row 125
$ws.Range("H125").Value = 4268
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4268
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 38412
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -43332

# Row 132
$ws.Range("H132").Value = 3539.1853
$ws.Range("I132").Value = 1728.2979
$ws.Range("J132").Value = 15698
$ws.Range("K132").Value = 5184.893700000001
$ws.Range("L132").Value = 47094
$ws.Range("M132").Value = -2654.893700000001
$ws.Range("N132").Value = -52154

# Row 135
$ws.Range("H135").Value = 17857802
$ws.Range("I135").Value = 527.97675
$ws.Range("J135").Value = 76924180
$ws.Range("K135").Value = 4751.79075
$ws.Range("L135").Value = 692317620
$ws.Range("M135").Value = -2216.79075
$ws.Range("N135").Value = -692322690

# Row 138
$ws.Range("H138").Value = 3392412.2
$ws.Range("I138").Value = 1136.9697
$ws.Range("J138").Value = 7696723
$ws.Range("K138").Value = 3410.9091
$ws.Range("L138").Value = 23090169
$ws.Range("M138").Value = 1729.0909
$ws.Range("N138").Value = -23100449


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3153.77
$ws.Range("I32").Value = 2681.0312
$ws.Range("J32").Value = 14499.5
$ws.Range("K32").Value = 2681.0312
$ws.Range("L32").Value = 14499.5
$ws.Range("M32").Value = -2394.0312
$ws.Range("N32").Value = -15073.5

# Row 88
$ws.Range("H88").Value = 2274.3
$ws.Range("I88").Value = 1909
$ws.Range("J88").Value = 2517.8333
$ws.Range("K88").Value = 1909
$ws.Range("L88").Value = 2517.8333
$ws.Range("M88").Value = -1503
$ws.Range("N88").Value = -3329.8333

# Row 91
$ws.Range("H91").Value = 2274.3
$ws.Range("I91").Value = 1909
$ws.Range("J91").Value = 2517.8333
$ws.Range("K91").Value = 1909
$ws.Range("L91").Value = 2517.8333
$ws.Range("M91").Value = -505
$ws.Range("N91").Value = -5325.8333


$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1942.9166
$ws.Range("I5").Value = 251
$ws.Range("J5").Value = 10402.5
$ws.Range("K5").Value = 251
$ws.Range("L5").Value = 10402.5
$ws.Range("M5").Value = -138
$ws.Range("N5").Value = -10628.5

# Row 86
$ws.Range("H86").Value = 1873.4736
$ws.Range("I86").Value = 1662.25
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1662.25
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -539.25
$ws.Range("N86").Value = -5246

# Row 89
$ws.Range("H89").Value = 1873.4736
$ws.Range("I89").Value = 1662.25
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 8311.25
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -2695.25
$ws.Range("N89").Value = -26232


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13890229
$ws.Range("I31").Value = 1249.1034
$ws.Range("J31").Value = 71430290
$ws.Range("K31").Value = 1249.1034
$ws.Range("L31").Value = 71430290
$ws.Range("M31").Value = -954.1034
$ws.Range("N31").Value = -71430880

# Row 34
$ws.Range("H34").Value = 13890229
$ws.Range("I34").Value = 1249.1034
$ws.Range("J34").Value = 71430290
$ws.Range("K34").Value = 1249.1034
$ws.Range("L34").Value = 71430290
$ws.Range("M34").Value = -1047.1034
$ws.Range("N34").Value = -71430694

# Row 50
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

# Row 58
$ws.Range("H58").Value = 2911.9575
$ws.Range("I58").Value = 3165.561
$ws.Range("J58").Value = 1179
$ws.Range("K58").Value = 3165.561
$ws.Range("L58").Value = 1179
$ws.Range("M58").Value = -2962.561
$ws.Range("N58").Value = -1585

# Row 132
$ws.Range("H132").Value = 825050.06
$ws.Range("I132").Value = 2075.282
$ws.Range("J132").Value = 6174386
$ws.Range("K132").Value = 6225.846
$ws.Range("L132").Value = 18523158
$ws.Range("M132").Value = -3695.846
$ws.Range("N132").Value = -18528218

# Row 136
$ws.Range("H136").Value = 2911.9575
$ws.Range("I136").Value = 3165.561
$ws.Range("J136").Value = 1179
$ws.Range("K136").Value = 9496.683
$ws.Range("L136").Value = 3537
$ws.Range("M136").Value = -6946.683000000001
$ws.Range("N136").Value = -8637


$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 3176.1365
$ws.Range("I58").Value = 1995
$ws.Range("K58").Value = 5985
$ws.Range("M58").Value = -5857

# Row 131
$ws.Range("H131").Value = 988.16
$ws.Range("I131").Value = 1001
$ws.Range("J131").Value = 988.0303
$ws.Range("K131").Value = 3003
$ws.Range("L131").Value = 2964.0909
$ws.Range("M131").Value = 2037
$ws.Range("N131").Value = -13044.0909


$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 9422
$ws.Range("I5").Value = 3312.5
$ws.Range("J5").Value = 13495
$ws.Range("K5").Value = 3312.5
$ws.Range("L5").Value = 13495
$ws.Range("M5").Value = -3200.5
$ws.Range("N5").Value = -13719


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1693.7273
$ws.Range("I16").Value = 1732.625
$ws.Range("J16").Value = 1590
$ws.Range("K16").Value = 1732.625
$ws.Range("L16").Value = 1590
$ws.Range("M16").Value = -1562.625
$ws.Range("N16").Value = -1930

# Row 82
$ws.Range("H82").Value = 1680
$ws.Range("I82").Value = 2700
$ws.Range("J82").Value = 1242.8572
$ws.Range("K82").Value = 2700
$ws.Range("L82").Value = 1242.8572
$ws.Range("M82").Value = -2339
$ws.Range("N82").Value = -1964.8572

# Row 85
$ws.Range("H85").Value = 1680
$ws.Range("I85").Value = 2700
$ws.Range("J85").Value = 1242.8572
$ws.Range("K85").Value = 2700
$ws.Range("L85").Value = 1242.8572
$ws.Range("M85").Value = -1452
$ws.Range("N85").Value = -3738.8572

# Row 132
$ws.Range("H132").Value = 3951.8164
$ws.Range("I132").Value = 3614.2563
$ws.Range("J132").Value = 5268.3
$ws.Range("K132").Value = 10842.7689
$ws.Range("L132").Value = 15804.9
$ws.Range("M132").Value = -8312.7689
$ws.Range("N132").Value = -20864.9


$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 78000
$ws.Range("I11").Value = 78000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 78000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -77858
$ws.Range("N11").ClearContents()

# Row 126
$ws.Range("H126").Value = 58824790
$ws.Range("I126").Value = 1210.75
$ws.Range("J126").Value = 111112424
$ws.Range("K126").Value = 3632.25
$ws.Range("L126").Value = 333337272
$ws.Range("M126").Value = -1162.25
$ws.Range("N126").Value = -333342212

# Row 132
$ws.Range("H132").Value = 4028.111
$ws.Range("I132").Value = 4636.6113
$ws.Range("J132").Value = 2811.111
$ws.Range("K132").Value = 13909.8339
$ws.Range("L132").Value = 8433.332999999999
$ws.Range("M132").Value = -11379.8339
$ws.Range("N132").Value = -13493.333

# Row 136
$ws.Range("H136").Value = 1625.9296
$ws.Range("I136").Value = 1420.7894
$ws.Range("J136").Value = 2461.1428
$ws.Range("K136").Value = 4262.3682
$ws.Range("L136").Value = 7383.428400000001
$ws.Range("M136").Value = -1712.3682
$ws.Range("N136").Value = -12483.4284

